$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells should match the style (bold, border, centered) of the
# existing header row, so copy formatting from the last existing header
# cell (AB1) before setting the new header text.
$ws.Range("AB1").Copy()
$ws.Range("AC1:AE1").PasteSpecial(-4122)

$ws.Range("AC1").Value = "Wins"
$ws.Range("AD1").Value = "Losses"
$ws.Range("AE1").Value = "Ties"

# Team record values for data rows 2-42
for ($r = 2; $r -le 42; $r++) {
    $ws.Range("AC$r").Value = 72
    $ws.Range("AD$r").Value = 90
    $ws.Range("AE$r").Value = 0
}
